# Weekly refresh of Fruta/Hortaliza "Terminal La Palmera de La Serena - Coco" data:
# each market-day row (2-41, except row 27 which is untouched) is re-stamped with
# the Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio ponderado /
# Precio $/Kg values pulled from that week's source row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns being refreshed, in order: Fecha(D), Volumen(M), Precio minimo(N),
# Precio maximo(O), Precio promedio ponderado(P), Precio $/Kg(S)
$columns = @("D", "M", "N", "O", "P", "S")

# Map of row number -> new [Fecha, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg]
$rowUpdates = @{
    2 = @(44315, 100, 20000, 21000, 20500, 1025)
    3 = @(44343, 100, 19500, 20000, 19750, 988)
    4 = @(44784, 160, 27000, 28000, 27500, 1375)
    5 = @(44778, 100, 23000, 24000, 23500, 1175)
    6 = @(44467, 200, 20000, 21000, 20500, 1025)
    7 = @(44365, 100, 20000, 21000, 20500, 1025)
    8 = @(44782, 200, 23500, 24000, 23750, 1188)
    9 = @(44410, 200, 20000, 21000, 20500, 1025)
    10 = @(44781, 160, 23000, 24000, 23500, 1175)
    11 = @(44448, 100, 20000, 21000, 20500, 1025)
    12 = @(44431, 160, 21000, 22000, 21500, 1075)
    13 = @(44333, 100, 19500, 20000, 19750, 988)
    14 = @(44434, 100, 20000, 21000, 20500, 1025)
    15 = @(44445, 160, 20000, 21000, 20500, 1025)
    16 = @(44441, 160, 20000, 21000, 20500, 1025)
    17 = @(44326, 160, 19500, 20000, 19750, 988)
    18 = @(44417, 160, 20000, 21000, 20500, 1025)
    19 = @(44420, 160, 20000, 21000, 20500, 1025)
    20 = @(44435, 260, 20000, 22000, 21115, 1056)
    21 = @(44879, 100, 28000, 30000, 29000, 1450)
    22 = @(44336, 100, 19500, 20000, 19750, 988)
    23 = @(44335, 200, 19000, 20000, 19500, 975)
    24 = @(44809, 60, 27000, 28000, 27500, 1375)
    25 = @(44418, 200, 20000, 21000, 20500, 1025)
    26 = @(44882, 120, 28000, 30000, 29000, 1450)
    28 = @(44776, 160, 23000, 24000, 23500, 1175)
    29 = @(44350, 160, 19000, 20000, 19500, 975)
    30 = @(44364, 140, 20000, 21000, 20500, 1025)
    31 = @(44407, 160, 20000, 21000, 20500, 1025)
    32 = @(44466, 100, 20000, 21000, 20500, 1025)
    33 = @(44427, 200, 20000, 21000, 20500, 1025)
    34 = @(44473, 40, 19500, 20000, 19750, 988)
    35 = @(44474, 200, 19000, 20000, 19500, 975)
    36 = @(44880, 100, 28000, 30000, 29000, 1450)
    37 = @(44428, 100, 20000, 21000, 20500, 1025)
    38 = @(44301, 100, 18000, 19000, 18500, 925)
    39 = @(44810, 100, 27000, 28000, 27500, 1375)
    40 = @(44442, 140, 20000, 21000, 20500, 1025)
    41 = @(44462, 100, 19500, 20000, 19750, 988)
}

foreach ($row in $rowUpdates.Keys) {
    $values = $rowUpdates[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}
